$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 14809.2
$ws.Range("I10").Value = 6962
$ws.Range("J10").Value = 26580
$ws.Range("K10").Value = 6962
$ws.Range("L10").Value = 26580
$ws.Range("M10").Value = -6669
$ws.Range("N10").Value = -27166

# Row 13
$ws.Range("H13").Value = 500
$ws.Range("J13").Value = 500
$ws.Range("L13").Value = 500
$ws.Range("N13").Value = -838

# Row 19
$ws.Range("H19").Value = 1234.1177
$ws.Range("J19").Value = 353.8
$ws.Range("L19").Value = 353.8
$ws.Range("N19").Value = -703.8

# Row 74
$ws.Range("H74").Value = 3308
$ws.Range("I74").Value = 1877.6
$ws.Range("K74").Value = 1877.6
$ws.Range("M74").Value = -941.5999999999999

# Row 77
$ws.Range("H77").Value = 3308
$ws.Range("I77").Value = 1877.6
$ws.Range("K77").Value = 9388
$ws.Range("M77").Value = -4708

# Row 92
$ws.Range("H92").Value = 12826.625
$ws.Range("I92").Value = 716.3333
$ws.Range("K92").Value = 716.3333
$ws.Range("M92").Value = 531.6667

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1950
$ws.Range("I2").Value = 1800
$ws.Range("J2").Value = 2100
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2100
$ws.Range("M2").Value = -1687
$ws.Range("N2").Value = -2326

# Row 32
$ws.Range("H32").Value = 2922.8293
$ws.Range("I32").Value = 2963.05
$ws.Range("K32").Value = 2963.05
$ws.Range("M32").Value = -2676.05

# Row 45
$ws.Range("H45").Value = 2757.1428

# Row 61
$ws.Range("H61").Value = 15155253
$ws.Range("I61").Value = 27781004
$ws.Range("J61").Value = 4351.4
$ws.Range("K61").Value = 27781004
$ws.Range("L61").Value = 4351.4
$ws.Range("M61").Value = -27780792
$ws.Range("N61").Value = -4775.4

# Row 110
$ws.Range("H110").Value = 111111630
$ws.Range("I110").Value = 111111630
$ws.Range("K110").Value = 111111630
$ws.Range("M110").Value = -111109585

# Row 116
$ws.Range("H116").Value = 1950
$ws.Range("I116").Value = 1800
$ws.Range("J116").Value = 2100
$ws.Range("K116").Value = 1800
$ws.Range("L116").Value = 2100
$ws.Range("M116").Value = 494
$ws.Range("N116").Value = -6688

# Row 136
$ws.Range("H136").Value = 15155253
$ws.Range("I136").Value = 27781004
$ws.Range("J136").Value = 4351.4
$ws.Range("K136").Value = 83343012
$ws.Range("L136").Value = 13054.2
$ws.Range("M136").Value = -83340462
$ws.Range("N136").Value = -18154.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1950
$ws.Range("I3").Value = 1800
$ws.Range("J3").Value = 2100
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2100
$ws.Range("M3").Value = -1686
$ws.Range("N3").Value = -2328

# Row 14
$ws.Range("H14").Value = 750
$ws.Range("J14").Value = 750
$ws.Range("L14").Value = 750
$ws.Range("N14").Value = -1094

$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 256.125
$ws.Range("I15").Value = 95
$ws.Range("K15").Value = 95
$ws.Range("M15").Value = 75

# Row 31
$ws.Range("H31").Value = 2626.862
$ws.Range("I31").Value = 1528.5278
$ws.Range("K31").Value = 1528.5278
$ws.Range("M31").Value = -1233.5278

# Row 34
$ws.Range("H34").Value = 2626.862
$ws.Range("I34").Value = 1528.5278
$ws.Range("K34").Value = 1528.5278
$ws.Range("M34").Value = -1326.5278

# Row 134
$ws.Range("H134").Value = 2605.75
$ws.Range("I134").Value = 1239.7693
$ws.Range("K134").Value = 3719.3079
$ws.Range("M134").Value = -1184.3079

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 1660.6
$ws.Range("J97").Value = 899
$ws.Range("L97").Value = 2697
$ws.Range("N97").Value = -3689

# Row 102
$ws.Range("H102").Value = 5500
$ws.Range("I102").Value = 7500
$ws.Range("J102").Value = 4833.3335
$ws.Range("K102").Value = 22500
$ws.Range("L102").Value = 14500.0005
$ws.Range("M102").Value = -20066
$ws.Range("N102").Value = -19368.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 499.5
$ws.Range("J13").Value = 499.5
$ws.Range("L13").Value = 499.5
$ws.Range("N13").Value = -777.5

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Row 22
$ws.Range("H22").Value = 2726.75
$ws.Range("I22").Value = 5004
$ws.Range("J22").Value = 449.5
$ws.Range("K22").Value = 5004
$ws.Range("L22").Value = 449.5
$ws.Range("M22").Value = -4475
$ws.Range("N22").Value = -1507.5

# Row 80
$ws.Range("H80").Value = 2813.4
$ws.Range("I80").Value = 2395
$ws.Range("K80").Value = 2395
$ws.Range("M80").Value = -1397

# Row 83
$ws.Range("H83").Value = 2813.4
$ws.Range("I83").Value = 2395
$ws.Range("K83").Value = 11975
$ws.Range("M83").Value = -6983

# Row 132
$ws.Range("H132").Value = 3412.6924
$ws.Range("I132").Value = 2610.6428
$ws.Range("K132").Value = 7831.928400000001
$ws.Range("M132").Value = -5301.928400000001

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 3666.3333
$ws.Range("I10").Value = 4500
$ws.Range("J10").Value = 1999
$ws.Range("K10").Value = 4500
$ws.Range("L10").Value = 1999
$ws.Range("M10").Value = -4360
$ws.Range("N10").Value = -2279

# Row 12
$ws.Range("H12").Value = 7859142.5
$ws.Range("I12").Value = 11000600
$ws.Range("J12").Value = 5498.5
$ws.Range("K12").Value = 11000600
$ws.Range("L12").Value = 5498.5
$ws.Range("M12").Value = -11000430
$ws.Range("N12").Value = -5838.5

# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Row 40
$ws.Range("H40").Value = 3311
$ws.Range("I40").Value = 2907.3333
$ws.Range("K40").Value = 2907.3333
$ws.Range("M40").Value = -2771.3333

# Row 82
$ws.Range("H82").Value = 912.86957
$ws.Range("I82").Value = 947
$ws.Range("J82").Value = 790
$ws.Range("K82").Value = 947
$ws.Range("L82").Value = 790
$ws.Range("M82").Value = -586
$ws.Range("N82").Value = -1512

# Row 85
$ws.Range("H85").Value = 912.86957
$ws.Range("I85").Value = 947
$ws.Range("J85").Value = 790
$ws.Range("K85").Value = 947
$ws.Range("L85").Value = 790
$ws.Range("M85").Value = 301
$ws.Range("N85").Value = -3286

# Row 132
$ws.Range("H132").Value = 6357.44
$ws.Range("I132").Value = 3772.7693
$ws.Range("K132").Value = 11318.3079
$ws.Range("M132").Value = -8788.3079

# Row 136
$ws.Range("H136").Value = 3042.6086
$ws.Range("I136").Value = 2117
$ws.Range("J136").Value = 5665.1665
$ws.Range("K136").Value = 6351
$ws.Range("L136").Value = 16995.4995
$ws.Range("M136").Value = -3801
$ws.Range("N136").Value = -22095.4995
